$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 24.02.2022 12:30"

# Update row 9 (OMV IKEA): D9 delta as a number, E9 as a real date/time value
$ws.Range("D9").Value = 0.4
$ws.Range("E9").Value = 44616.51069444444
$ws.Range("E9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
